{"js": "// Split the final paragraph (\"Hi Khanh Huyen\") into two paragraphs by\n// inserting a new paragraph after it containing \"Xin ch\u00e0o t\u1ea5t c\u1ea3 m\u1ecdi ng\u01b0\u1eddi\",\n// then relocate the trailing _GoBack bookmark from the end of the first\n// paragraph to the end of the newly added second paragraph (matching how\n// Word leaves _GoBack at the most recent edit point).\n\nconst doc = context.document;\nconst body = doc.body;\n\n// 1. Grab the existing (only) paragraph and add the new sentence right\n//    after it. insertParagraph inherits the paragraph/run formatting\n//    (lang=\"vi-VN\") from the paragraph it is inserted relative to.\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items\");\nawait context.sync();\n\nconst firstParagraph = paragraphs.items[0];\nfirstParagraph.insertParagraph(\"Xin ch\u00e0o t\u1ea5t c\u1ea3 m\u1ecdi ng\u01b0\u1eddi\", Word.InsertLocation.after);\nawait context.sync();\n\n// 2. Re-fetch the paragraphs fresh from the body instead of reusing the\n//    paragraph handle that insertParagraph returned \u2014 the freshly-queried\n//    object gives a reliable collapsed range at its end.\nconst refreshedParagraphs = body.paragraphs;\nrefreshedParagraphs.load(\"items\");\nawait context.sync();\nconst newParagraph = refreshedParagraphs.items[refreshedParagraphs.items.length - 1];\n\n// 3. Move the _GoBack bookmark so it collapses right after the new text,\n//    matching where Word would leave it after the most recent typing.\ndoc.deleteBookmark(\"_GoBack\");\nconst endOfNewParagraph = newParagraph.getRange(Word.RangeLocation.end);\nendOfNewParagraph.insertBookmark(\"_GoBack\");\n\nawait context.sync();\n", "ps1": "# Split the final paragraph (\"Hi Khanh Huyen\") into two paragraphs by\n# inserting a new paragraph after it containing \"Xin ch\u00e0o t\u1ea5t c\u1ea3 m\u1ecdi ng\u01b0\u1eddi\",\n# then relocate the trailing _GoBack bookmark from the end of the first\n# paragraph to the end of the newly added second paragraph (matching how\n# Word leaves _GoBack at the most recent edit point).\n\n$d = $word.ActiveDocument\n\n# 1. Insert a new paragraph right after the existing \"Hi Khanh Huyen\" paragraph,\n#    inheriting its paragraph/run formatting (lang=\"vi-VN\").\n$firstPara = $d.Paragraphs.First\n$firstPara.Range.InsertParagraphAfter()\n\n# 2. Fill the new (now-last) paragraph with the new sentence. A trailing\n#    sentinel character is appended first so the real insertion point for the\n#    bookmark below is never the very last position(s) of the document\n#    (collapsed ranges flush with the document end are not positioned\n#    reliably), then the sentinel is stripped back out afterward.\n$newPara = $d.Paragraphs.Last\n$newPara.Range.Text = \"Xin ch\u00e0o t\u1ea5t c\u1ea3 m\u1ecdi ng\u01b0\u1eddi\" + \"X\"\n\n# 3. Move the _GoBack bookmark so it collapses right after the new text\n#    (before the sentinel / paragraph mark), matching where Word would leave\n#    it after the most recent typing.\n$goBack = $d.Bookmarks(\"_GoBack\")\n$goBack.Delete()\n\n$anchorPos = $d.Content.End - 2\n$anchorRange = $d.Range($anchorPos, $anchorPos)\n$d.Bookmarks.Add(\"_GoBack\", $anchorRange)\n\n# 4. Remove the sentinel character now that the bookmark is anchored.\n$sentinelRange = $d.Range($d.Content.End - 2, $d.Content.End - 1)\n$sentinelRange.Delete()\n"}
